$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the oldest data row (old row 2, date 2007-11-14 / serial 39400).
# This shifts every subsequent row up by one and also shrinks the used
# range from A1:E19 down to A1:E18.
$ws.Rows("2").Delete()

# Recalculated y_1_forecast values (column E) for the remaining rows.
$e = @{
  2  = 1.560682679516057
  3  = 2.1453644888767
  4  = 1.317672174811868
  5  = 1.501816644427989
  6  = 1.028888107831327
  7  = 1.303605130836716
  8  = 1.192378712846454
  9  = 1.210961441871872
  10 = 2.033218171624651
  11 = 2.152537330144288
  12 = 1.966855307908655
  13 = 1.950353221540246
  14 = 2.210985773414453
  15 = 1.114171399050901
  16 = 0.1338254721205745
  17 = 1.823564868738359
  18 = 0.6266145540918089
}

foreach ($r in $e.Keys) {
  $ws.Cells.Item($r, 5).Value = $e[$r]
}
